$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (old row 6) so the table shrinks from 6 rows to 5
$ws.Rows.Item(6).Delete()

# Widen the columns that grew from width 7 to width 8
# (ColumnWidth uses character units; 7.16666... round-trips to a stored width of exactly 8)
$wideColIndexes = @(2,3,7,10,11,12,13,15,17,22,24,27,28,29,30,34)
foreach ($idx in $wideColIndexes) {
    $ws.Columns.Item($idx).ColumnWidth = 7.166666666666667
}

# Update data values for rows 2-5 (new sensor readings)
# Row 2
$ws.Cells.Item(2, 1).Value = 45179.50694444445
$ws.Cells.Item(2, 2).Value = 7.205
$ws.Cells.Item(2, 3).Value = 5.605
$ws.Cells.Item(2, 4).Value = 1.192
$ws.Cells.Item(2, 5).Value = 15.61
$ws.Cells.Item(2, 6).Value = 12.374
$ws.Cells.Item(2, 7).Value = 4.794
$ws.Cells.Item(2, 8).Value = 14.798
$ws.Cells.Item(2, 9).Value = 8.957000000000001
$ws.Cells.Item(2, 10).Value = 4.429
$ws.Cells.Item(2, 11).Value = 5.631
$ws.Cells.Item(2, 12).Value = 6.249
$ws.Cells.Item(2, 13).Value = 7.306
$ws.Cells.Item(2, 14).Value = 2.788
$ws.Cells.Item(2, 15).Value = 6.015
$ws.Cells.Item(2, 16).Value = 7.738
$ws.Cells.Item(2, 17).Value = 5.138
$ws.Cells.Item(2, 18).Value = 0.492
$ws.Cells.Item(2, 19).Value = 0.931
$ws.Cells.Item(2, 20).Value = 84.14
$ws.Cells.Item(2, 21).Value = 16.424
$ws.Cells.Item(2, 22).Value = 5.552
$ws.Cells.Item(2, 23).Value = 10.174
$ws.Cells.Item(2, 24).Value = 6.283
$ws.Cells.Item(2, 25).Value = 0.894
$ws.Cells.Item(2, 26).Value = 9.741
$ws.Cells.Item(2, 27).Value = 4.361
$ws.Cells.Item(2, 28).Value = 4.885
$ws.Cells.Item(2, 29).Value = 6.06
$ws.Cells.Item(2, 30).Value = 8.134
$ws.Cells.Item(2, 31).Value = 1.522
$ws.Cells.Item(2, 32).Value = 13.18
$ws.Cells.Item(2, 33).Value = 3.625
$ws.Cells.Item(2, 34).Value = 6.387

# Row 3
$ws.Cells.Item(3, 1).Value = 45179.51388888889
$ws.Cells.Item(3, 2).Value = 17.977
$ws.Cells.Item(3, 3).Value = 13.558
$ws.Cells.Item(3, 4).Value = 1.001
$ws.Cells.Item(3, 5).Value = 39.344
$ws.Cells.Item(3, 6).Value = 32.009
$ws.Cells.Item(3, 7).Value = 13.606
$ws.Cells.Item(3, 8).Value = 50.673
$ws.Cells.Item(3, 9).Value = 21.876
$ws.Cells.Item(3, 10).Value = 10.256
$ws.Cells.Item(3, 11).Value = 14.364
$ws.Cells.Item(3, 12).Value = 15.732
$ws.Cells.Item(3, 13).Value = 17.031
$ws.Cells.Item(3, 14).Value = 5.112
$ws.Cells.Item(3, 15).Value = 14.287
$ws.Cells.Item(3, 16).Value = 19.947
$ws.Cells.Item(3, 17).Value = 12.037
$ws.Cells.Item(3, 18).Value = 0.3
$ws.Cells.Item(3, 19).Value = 0.887
$ws.Cells.Item(3, 20).Value = 210.126
$ws.Cells.Item(3, 21).Value = 39.874
$ws.Cells.Item(3, 22).Value = 13.187
$ws.Cells.Item(3, 23).Value = 26.488
$ws.Cells.Item(3, 24).Value = 14.42
$ws.Cells.Item(3, 25).Value = 1.903
$ws.Cells.Item(3, 26).Value = 26.733
$ws.Cells.Item(3, 27).Value = 11.349
$ws.Cells.Item(3, 28).Value = 10.621
$ws.Cells.Item(3, 29).Value = 12.578
$ws.Cells.Item(3, 30).Value = 17.239
$ws.Cells.Item(3, 31).Value = 0.707
$ws.Cells.Item(3, 32).Value = 46.499
$ws.Cells.Item(3, 33).Value = 7.797
$ws.Cells.Item(3, 34).Value = 16.184

# Row 4
$ws.Cells.Item(4, 1).Value = 45179.52083333334
$ws.Cells.Item(4, 2).Value = 22.383
$ws.Cells.Item(4, 3).Value = 16.832
$ws.Cells.Item(4, 4).Value = 1.002
$ws.Cells.Item(4, 5).Value = 48.928
$ws.Cells.Item(4, 6).Value = 40.032
$ws.Cells.Item(4, 7).Value = 17.221
$ws.Cells.Item(4, 8).Value = 67.917
$ws.Cells.Item(4, 9).Value = 27.171
$ws.Cells.Item(4, 10).Value = 12.567
$ws.Cells.Item(4, 11).Value = 17.962
$ws.Cells.Item(4, 12).Value = 19.59
$ws.Cells.Item(4, 13).Value = 20.992
$ws.Cells.Item(4, 14).Value = 6.052
$ws.Cells.Item(4, 15).Value = 17.67
$ws.Cells.Item(4, 16).Value = 24.934
$ws.Cells.Item(4, 17).Value = 14.813
$ws.Cells.Item(4, 18).Value = 0.265
$ws.Cells.Item(4, 19).Value = 0.896
$ws.Cells.Item(4, 20).Value = 261.648
$ws.Cells.Item(4, 21).Value = 49.415
$ws.Cells.Item(4, 22).Value = 16.31
$ws.Cells.Item(4, 23).Value = 33.105
$ws.Cells.Item(4, 24).Value = 17.733
$ws.Cells.Item(4, 25).Value = 2.331
$ws.Cells.Item(4, 26).Value = 34.098
$ws.Cells.Item(4, 27).Value = 14.201
$ws.Cells.Item(4, 28).Value = 12.945
$ws.Cells.Item(4, 29).Value = 15.26
$ws.Cells.Item(4, 30).Value = 21.024
$ws.Cells.Item(4, 31).Value = 0.461
$ws.Cells.Item(4, 32).Value = 62.08
$ws.Cells.Item(4, 33).Value = 9.497999999999999
$ws.Cells.Item(4, 34).Value = 20.184

# Row 5
$ws.Cells.Item(5, 1).Value = 45179.52777777778
$ws.Cells.Item(5, 2).Value = 5.13
$ws.Cells.Item(5, 3).Value = 3.86
$ws.Cells.Item(5, 4).Value = 0.36
$ws.Cells.Item(5, 5).Value = 11.38
$ws.Cells.Item(5, 6).Value = 9.06
$ws.Cells.Item(5, 7).Value = 3.73
$ws.Cells.Item(5, 8).Value = 22.12
$ws.Cells.Item(5, 9).Value = 6.26
$ws.Cells.Item(5, 10).Value = 3.17
$ws.Cells.Item(5, 11).Value = 3.97
$ws.Cells.Item(5, 12).Value = 4.53
$ws.Cells.Item(5, 13).Value = 5
$ws.Cells.Item(5, 14).Value = 1.63
$ws.Cells.Item(5, 15).Value = 4.14
$ws.Cells.Item(5, 16).Value = 5.73
$ws.Cells.Item(5, 17).Value = 3.58
$ws.Cells.Item(5, 18).Value = 0.09
$ws.Cells.Item(5, 19).Value = 0.33
$ws.Cells.Item(5, 20).Value = 55.68
$ws.Cells.Item(5, 21).Value = 11.88
$ws.Cells.Item(5, 22).Value = 3.82
$ws.Cells.Item(5, 23).Value = 7.81
$ws.Cells.Item(5, 24).Value = 4.29
$ws.Cells.Item(5, 25).Value = 0.55
$ws.Cells.Item(5, 26).Value = 10.76
$ws.Cells.Item(5, 27).Value = 3.21
$ws.Cells.Item(5, 28).Value = 3.21
$ws.Cells.Item(5, 29).Value = 3.78
$ws.Cells.Item(5, 30).Value = 5.06
$ws.Cells.Item(5, 31).Value = 0.34
$ws.Cells.Item(5, 32).Value = 20.8
$ws.Cells.Item(5, 33).Value = 2.33
$ws.Cells.Item(5, 34).Value = 4.62

